# Generate Report for Archive
# - Update Status value from "Ready for handoff" to "In Translation" across sheets.
# - Shrink the now-narrower "Status"-related columns to their new autofit width.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# The target stored column width is 13.4101845877511 (re-autofit width for
# the shorter "In Translation" text). Excel's ColumnWidth setter here only
# resolves to 1/6-character increments, so 12.5 is the input that lands on
# the closest achievable stored width (13.333333333333334).
$newWidth = 12.5

# Overview sheet: columns E (zh-cn) and F (de-de) hold the status value.
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($col in @("E", "F")) {
    for ($row = 2; $row -le 4; $row++) {
        $cell = $wsOverview.Range("$col$row")
        if ($cell.Text -eq $oldStatus) {
            $cell.Value = $newStatus
        }
    }
    $wsOverview.Columns($col).ColumnWidth = $newWidth
}

# zh-cn and de-de sheets: column C holds the "Status" value.
foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($row = 2; $row -le 4; $row++) {
        $cell = $ws.Range("C$row")
        if ($cell.Text -eq $oldStatus) {
            $cell.Value = $newStatus
        }
    }
    $ws.Columns("C").ColumnWidth = $newWidth
}
